# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (E) / "Valor Mora" (F) pair for the two periods belonging
# to worker CC 9186025 (JAVIER CAJAR CAJAR NAVARRO) are swapped between rows
# 16 and 17, and the "Salario Basico" (G) figure is refreshed for every row
# in the table (rows 16-18) to reflect the updated account-statement data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: period 1804 (2951 / 781242) -> period 1701 (27578 / 737717)
$ws.Range("E16").Value = "1701"
$ws.Range("F16").Value = 27578
$ws.Range("G16").Value = 737717

# Row 17: period 1701 (27578 / 781242) -> period 1804 (2951 / 737717)
$ws.Range("E17").Value = "1804"
$ws.Range("F17").Value = 2951
$ws.Range("G17").Value = 737717

# Row 18: period stays 1804 (2951), only the Salario Basico is refreshed
$ws.Range("G18").Value = 737717
